$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Apply same style as the other header cells (A1:E1) to the new headers
$ws.Range("A1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill F2:H21 with boolean FALSE values
$ws.Range("F2:H21").Value = $false
